$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: Target cluster FAPs -> ECs ---
$ws.Range("D2").Value = "ECs"
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2956425
$ws.Range("N2").Value = 0.5912849999999999
$ws.Range("O2").Value = 0.179122830242808
$ws.Range("P2").Value = 0.1717318139150154
$ws.Range("Q2").Value = 0.004596551042499999
$ws.Range("R2").Value = 0.02757930625499999
$ws.Range("S2").Value = 0.179122830242808
$ws.Range("T2").Value = 0.1717318139150154

# --- Row 3: Target cluster MuSCs -> FAPs ---
$ws.Range("D3").Value = "FAPs"
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.142069
$ws.Range("N3").Value = 0.426207
$ws.Range("O3").Value = 0.08607626227543567
$ws.Range("P3").Value = 0.1237868391947656
$ws.Range("Q3").Value = 0.002208841455666667
$ws.Range("R3").Value = 0.019879573101
$ws.Range("S3").Value = 0.08607626227543567
$ws.Range("T3").Value = 0.1237868391947656

# --- Row 4: Target cluster Neutrophils -> MuSCs ---
$ws.Range("D4").Value = "MuSCs"
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.21279
$ws.Range("N4").Value = 2.42558
$ws.Range("O4").Value = 0.7348009074817563
$ws.Range("P4").Value = 0.704481346890219
$ws.Range("Q4").Value = 0.01885605465666667
$ws.Range("R4").Value = 0.11313632794
$ws.Range("S4").Value = 0.7348009074817563
$ws.Range("T4").Value = 0.704481346890219
